# Update "still cancelled" figures on the 3rd tab / summary worksheet.
# Diff only touches raw-value cells (Current Month Active [C], Lifetime
# Cancels [D], Residual Core Plans [J], Current Active DVH [M]); all
# dependent formula cells (E,H,I,K,O,P,T,U) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = 974
$ws.Range("M3").Value  = 108

$ws.Range("M5").Value  = 54

$ws.Range("C6").Value  = 151
$ws.Range("J6").Value  = 1

$ws.Range("C7").Value  = 951
$ws.Range("M7").Value  = 238

$ws.Range("C8").Value  = 937

$ws.Range("C9").Value  = 389
$ws.Range("M9").Value  = 90

$ws.Range("M10").Value = 79

$ws.Range("C12").Value = 397

$ws.Range("C13").Value = 111
$ws.Range("J13").Value = 1

$ws.Range("C14").Value = 155
$ws.Range("M14").Value = 15

$ws.Range("M16").Value = 112

$ws.Range("C17").Value = 591

$ws.Range("C19").Value = 626
$ws.Range("J19").Value = 56
$ws.Range("M19").Value = 103

$ws.Range("M20").Value = 78

$ws.Range("C21").Value = 1141
$ws.Range("J21").Value = 516
$ws.Range("M21").Value = 98

$ws.Range("D22").Value = 2
$ws.Range("M22").Value = 59

$ws.Range("C23").Value = 427
$ws.Range("J23").Value = 123
$ws.Range("M23").Value = 101

$ws.Range("C24").Value = 389
$ws.Range("J24").Value = 1
$ws.Range("M24").Value = 84

$wb.Save()
